$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null
}

# Italicize the $occurrence-th (1-based) appearance of $phrase inside a given table cell.
function Italicize-InCell($tableIndex, $rowIndex, $colIndex, $phrase, $occurrence) {
    $cellRange = $d.Tables.Item($tableIndex).Cell($rowIndex, $colIndex).Range
    $cellStart = $cellRange.Start
    $text = $cellRange.Text
    $searchFrom = 0
    $idx = -1
    for ($i = 1; $i -le $occurrence; $i++) {
        $idx = $text.IndexOf($phrase, $searchFrom)
        if ($idx -lt 0) { break }
        $searchFrom = $idx + 1
    }
    if ($idx -lt 0) {
        Write-Output "PHRASE NOT FOUND: $phrase"
        return
    }
    $wordStart = $cellStart + $idx
    $wordEnd = $wordStart + $phrase.Length
    $wordRange = $d.Range($wordStart, $wordEnd)
    $wordRange.Italic = 1
}

# ---------------------------------------------------------------------------
# PA "Vartotojo registracija" (Table 1) - registration
# ---------------------------------------------------------------------------
Replace-Text "Vartotojas atsidarė registravimo langą." "Vartotojas turi būti neprisiregistravęs prie sistemos."
Replace-Text "Vartotojas nori užsiregistruoti sistemoje." "Vartotojas atsidarė registravimo langą."

# ---------------------------------------------------------------------------
# PA "Vartotojo prisijungimas" (Table 2) - login
# ---------------------------------------------------------------------------
Replace-Text "Vartotojas nori prisijungti prie sistemos." "Vartotojas atsidarė prisijungimo langą."

# ---------------------------------------------------------------------------
# PA "Profilio administravimas" (Table 3) - profile
# ---------------------------------------------------------------------------
Replace-Text "Vartotojas turi būti prisijungęs prie sistemos ir atsidaręs savo profilio langą." "Vartotojas turi būti prisijungęs prie sistemos."
Replace-Text "Vartotojas nori pasikeisti savo duomenis arba pašalinti paskyrą." "Vartotojas atsidarė savo profilio langą."
Italicize-InCell 3 6 2 "profilio" 1

# ---------------------------------------------------------------------------
# PA "Filmų vertinimas/komentavimas" (Table 4) - rating/commenting
# ---------------------------------------------------------------------------
Replace-Text "Vartotojas turi būti prisijungęs prie sistemos kaip klientas ir atsidaręs filmo langą." "Vartotojas turi būti prisijungęs prie sistemos kaip klientas."
Replace-Text "Vartotojas nori įvertinti ir/arba pakomentuoti norimą filmą." "Vartotojas atsidarė filmo langą ir įvertina ir/arba pakomentuoja filmą."
Italicize-InCell 4 6 2 "filmo" 1

# ---------------------------------------------------------------------------
# PA "Dalyvavimas balsavimuose" (Table 5) - voting
# ---------------------------------------------------------------------------
Replace-Text "Vartotojas turi būti prisijungęs prie sistemos kaip klientas ir atsidaręs balsavimo langą." "Vartotojas turi būti prisijungęs prie sistemos kaip klientas."
Replace-Text "Vartotojas nori balsuoti už norimą kino kūrėją." "Vartotojas atsidarė balsavimo langą ir balsuoja už norimą kino kūrėją."
Italicize-InCell 5 6 2 "balsavimo" 1

# ---------------------------------------------------------------------------
# PA "Kino teatrų prenumeravimas" (Table 6) - cinema subscription
# ---------------------------------------------------------------------------
Replace-Text "Vartotojas turi būti prisijungęs prie sistemos kaip klientas ir atsidaręs kino teatro langą." "Vartotojas turi būti prisijungęs prie sistemos kaip klientas."
Replace-Text "Vartotojas nori prenumeruoti kino teatrą." "Vartotojas atsidarė kino teatro langą ir prenumeruoja kino teatrą."
Italicize-InCell 6 6 2 "kino teatro" 1

# ---------------------------------------------------------------------------
# PA "Bilieto rezervavimas" (Table 7) - ticket reservation
# ---------------------------------------------------------------------------
Replace-Text "Vartotojas turi būti prisijungęs prie sistemos kaip klientas ir atsidaręs kino teatro filmo langą." "Vartotojas turi būti prisijungęs prie sistemos kaip klientas."
Replace-Text "Vartotojas nori rezervuoti vietą į kino filmo seansą." "Vartotojas atsidarė kino teatro filmo langą ir rezervuoja vietą į kino filmo seansą."
Italicize-InCell 7 6 2 "kino teatro filmo" 1

# ---------------------------------------------------------------------------
# Move the _GoBack bookmark to reflect the position of the final edit
# (between "rezervuoja" and " vietą į kino filmo seansą.")
# ---------------------------------------------------------------------------
$bm = $d.Bookmarks
if ($bm.Exists("_GoBack")) {
    $bm.Item("_GoBack").Delete() | Out-Null
}
$lastCell = $d.Tables.Item(7).Cell(6, 2).Range
$lastCellStart = $lastCell.Start
$lastText = $lastCell.Text
$marker = "rezervuoja"
$markerIdx = $lastText.IndexOf($marker)
if ($markerIdx -ge 0) {
    $pos = $lastCellStart + $markerIdx + $marker.Length
    $bmRange = $d.Range($pos, $pos)
    $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
}

Write-Output "Done"
